$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing "weekly_activity" rows' scores (columns C = "Friss mi Stoub", D = "Smash di weg")
$ws.Range("C13").Value = $true
$ws.Range("C16").Value = $true
$ws.Range("D19").Value = $true
$ws.Range("D22").Value = $true
$ws.Range("C34").Value = $true
$ws.Range("C37").Value = $true
$ws.Range("D37").Value = $true
$ws.Range("D40").Value = $true
$ws.Range("D43").Value = $true
$ws.Range("C55").Value = $true
$ws.Range("C58").Value = $true
$ws.Range("D58").Value = $true
$ws.Range("C61").Value = $true
$ws.Range("D61").Value = $true

# Append three new days (2025-02-21, 2025-02-22, 2025-02-23) of daily motivation scores
$newRows = @(
    @("2025-02-21", "sleep",           $false, $true),
    @("2025-02-21", "activity",        $true,  $true),
    @("2025-02-21", "weekly_activity", $true,  $true),
    @("2025-02-22", "sleep",           $true,  $false),
    @("2025-02-22", "activity",        $true,  $true),
    @("2025-02-22", "weekly_activity", $false, $false),
    @("2025-02-23", "sleep",           $true,  $false),
    @("2025-02-23", "activity",        $false, $false),
    @("2025-02-23", "weekly_activity", $false, $false)
)

$r = 62
foreach ($row in $newRows) {
    $dateCell = $ws.Cells.Item($r, 1)
    $dateCell.Value = "'" + $row[0]
    $dateCell.Style = "Normal"

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
